$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A new harbour ("las Palmas") is added as the first port-of-call; the
# previous harbour1 value ("Bremerhaven") shifts over to the next column,
# and the stale "Cape Town" value that used to sit there is replaced.
$ws.Range("D2").Value = "las Palmas"
$ws.Range("E2").Value = "Bremerhaven"

# Reflect the cell range the author had selected when the file was saved.
$ws.Range("B4:C9").Select() | Out-Null
